$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: apply a batch of F-column (numeric) updates to a worksheet, plus
# any I-column (cover image URL) text updates.
# ---------------------------------------------------------------------------
function Set-FValues {
    param($ws, $updates)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) — F-column "want to go" counters bumped, plus a
# refreshed cover image for row 19.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$expoUpdates = @{
    5  = 1702
    6  = 3286
    7  = 922
    8  = 2112
    9  = 2030
    10 = 1047
    11 = 563
    13 = 1634
    14 = 354
    17 = 78
    18 = 113
    19 = 1483
    20 = 556
    21 = 658
    22 = 343
    23 = 11885
    24 = 11911
    25 = 873
    26 = 671
    27 = 11
    29 = 167
    30 = 483
}
Set-FValues $wsExpo $expoUpdates
$wsExpo.Cells.Item(19, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/rrEX7kZn1715245404853.jpeg"

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances) — insert a new show as row 7 ("广州·跨越二次元
# ACG神级动漫世界巡回演唱会"), pushing the old row 7 ("一个陌生女人的来信")
# down to row 8.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Rows.Item(7).Insert()

$wsShow.Cells.Item(7, 1).Value = 6
$wsShow.Cells.Item(7, 1).Font.Bold = $true
$wsShow.Cells.Item(7, 1).HorizontalAlignment = -4108
$wsShow.Cells.Item(7, 1).VerticalAlignment = -4160
$wsShow.Cells.Item(7, 1).Borders.LineStyle = 1
$wsShow.Cells.Item(7, 2).Value = "'2024-07-20"
$wsShow.Cells.Item(7, 3).Value = "广州·跨越二次元ACG神级动漫世界巡回演唱会"
$wsShow.Cells.Item(7, 4).Value = "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院"
$wsShow.Cells.Item(7, 5).Value = "2024.07.20 19:30-07.20 21:10"
$wsShow.Cells.Item(7, 6).Value = 1
$wsShow.Cells.Item(7, 7).Value = 120
$wsShow.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85353"
$wsShow.Cells.Item(7, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg"

# The pushed-down row keeps its own data but its running index (col A) bumps.
$wsShow.Cells.Item(8, 1).Value = 7

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) — same F-column bumps as 展览 (at shifted row
# numbers) plus the same new 演出 row inserted before "一个陌生女人的来信".
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$allUpdates = @{
    7  = 1702
    8  = 3286
    9  = 922
    10 = 2112
    11 = 2030
    12 = 1047
    13 = 563
    15 = 1634
    16 = 354
    20 = 78
    22 = 113
    23 = 1483
    24 = 556
    25 = 658
    26 = 343
    27 = 11885
    28 = 11911
    29 = 873
    30 = 671
    31 = 11
    35 = 167
    36 = 483
}
Set-FValues $wsAll $allUpdates
$wsAll.Cells.Item(23, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/rrEX7kZn1715245404853.jpeg"

$wsAll.Rows.Item(37).Insert()

$wsAll.Cells.Item(37, 1).Value = 36
$wsAll.Cells.Item(37, 1).Font.Bold = $true
$wsAll.Cells.Item(37, 1).HorizontalAlignment = -4108
$wsAll.Cells.Item(37, 1).VerticalAlignment = -4160
$wsAll.Cells.Item(37, 1).Borders.LineStyle = 1
$wsAll.Cells.Item(37, 2).Value = "'2024-07-20"
$wsAll.Cells.Item(37, 3).Value = "广州·跨越二次元ACG神级动漫世界巡回演唱会"
$wsAll.Cells.Item(37, 4).Value = "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院"
$wsAll.Cells.Item(37, 5).Value = "2024.07.20 19:30-07.20 21:10"
$wsAll.Cells.Item(37, 6).Value = 1
$wsAll.Cells.Item(37, 7).Value = 120
$wsAll.Cells.Item(37, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85353"
$wsAll.Cells.Item(37, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg"

$wsAll.Cells.Item(38, 1).Value = 37

Write-Output "edit complete"
